# Update API Posture Management workbook:
#  - Legend sheet gains a header row ("Column1"/"Column2") and becomes an
#    Excel Table (Table1)
#  - Sheet1/Sheet2 are renamed to Data/Legend
#  - Legend becomes the active (selected) sheet/tab

$wb = $excel.ActiveWorkbook

$wsData   = $wb.Worksheets.Item(1)
$wsLegend = $wb.Worksheets.Item(2)

# --- Legend sheet: insert a header row above the existing lookup table ---
$wsLegend.Rows.Item(1).Insert()
$wsLegend.Range("A1").Value = "Column1"
$wsLegend.Range("B1").Value = "Column2"

# Turn A1:B7 into a real Excel Table (ListObject), headers already in place
$table = $wsLegend.ListObjects.Add(1, $wsLegend.Range("A1:B7"), $null, 1)
$table.Name = "Table1"

# --- Rename sheets ---
$wsData.Name = "Data"
$wsLegend.Name = "Legend"

# --- Make Legend the active/selected sheet ---
$wsLegend.Range("A1:B7").Select() | Out-Null
$wsLegend.Activate()
